# Automatyczna aktualizacja Excela [2025-07-27 06:26:26]
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "powiat krakowski"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Rows whose "Data ostatniej aktualizacji" (column E) moves from 2025-07-24 to 2025-07-27
$ws1_E_rows = @(2,3,4,5,7,9,10,11,12,14,16,17,18,22,23,24,25,26,27,30,31,34,35,37,38,39,41,42,43,44,46,48,50,51,52,53,54,55,57,58,59,60,61,62,83,84,85,86,87,88,89,90,91,92,93,94,96,97)
foreach ($r in $ws1_E_rows) {
    $cell = $ws1.Cells.Item($r, 5)
    $cell.NumberFormat = "@"
    $cell.Value = "2025-07-27"
}

# Rows whose "Aktywne" flag (column H) flips from TRUE to FALSE
$ws1_H_rows = @(6,13,36,65,66,68,78,82,95,98,99,100,101,102,103,104,105,106,107,108,109,110)
foreach ($r in $ws1_H_rows) {
    $ws1.Cells.Item($r, 8).Value = $false
}

# ---------------------------------------------------------------------------
# Sheet 2: "powiat wielicki"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Rows whose "Data ostatniej aktualizacji" (column E) moves from 2025-07-24 to 2025-07-27
$ws2_E_rows = @(2,3,4,5,6,7,8,9,13,17,23,24,26,29,30,33,34,35,36,37,38,39,40,42,43,44,46,47,48,49,50,51,52,62,63,65,66,68,69,70,71,72,73,74,75,78,79)
foreach ($r in $ws2_E_rows) {
    $cell = $ws2.Cells.Item($r, 5)
    $cell.NumberFormat = "@"
    $cell.Value = "2025-07-27"
}

# Rows whose "Aktywne" flag (column H) flips from TRUE to FALSE
$ws2_H_rows = @(31,56,77,80)
foreach ($r in $ws2_H_rows) {
    $ws2.Cells.Item($r, 8).Value = $false
}

# New listings appended at the bottom of sheet 2 (rows 81-85)
$newRows = @(
    @("Działka Budowlana | 16 ar | Jawczyce", "Jawczyce, Biskupice, wielicki, małopolskie", 99000, "2025-07-27", "2025-07-27", 99000, 0.49, $true, "https://www.otodom.pl/pl/oferta/dzialka-budowlana-16-ar-jawczyce-ID4xr4L"),
    @("Działka budowlana | Koźmice Wielkie | gm. Wieliczka | 1416 m²", "Koźmice Wielkie, Wieliczka, wielicki, małopolskie", 249000, "2025-07-27", "2025-07-27", 249000, 0.49, $true, "https://www.otodom.pl/pl/oferta/dzialka-budowlana-kozmice-wielkie-gm-wieliczka-1416-m-ID4xqLp"),
    @("14a budowlane przy drodze z mediami, Krk Wieliczka", "Sygneczów, Wieliczka, wielicki, małopolskie", 152500, "2025-07-27", "2025-07-27", 152500, 0.49, $true, "https://www.otodom.pl/pl/oferta/14a-budowlane-przy-drodze-z-mediami-krk-wieliczka-ID4xpge"),
    @("Kobylec | Działka budowlana przy Lesie -  ok. 14 ar / 179.000 zł", "Kobylec, Łapanów, bocheński, małopolskie", 179000, "2025-07-27", "2025-07-27", 179000, 0.49, $true, "https://www.otodom.pl/pl/oferta/kobylec-dzialka-budowlana-przy-lesie-ok-14-ar-179-000-zl-ID4xjra"),
    @("Działka Budowlana | 16 ar | Jawczyce", "Jawczyce, Biskupice, wielicki, małopolskie", 99000, "2025-07-27", "2025-07-27", 99000, 0.49, $true, "https://www.otodom.pl/hpr/pl/oferta/dzialka-budowlana-16-ar-jawczyce-ID4xr4L")
)

$startRow = 81
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws2.Cells.Item($r, 1).Value = $data[0]
    $ws2.Cells.Item($r, 2).Value = $data[1]
    $ws2.Cells.Item($r, 3).Value = $data[2]

    $dCell = $ws2.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $data[3]

    $eCell = $ws2.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $data[4]

    $ws2.Cells.Item($r, 6).Value = $data[5]
    $ws2.Cells.Item($r, 7).Value = $data[6]
    $ws2.Cells.Item($r, 8).Value = $data[7]
    $ws2.Cells.Item($r, 9).Value = $data[8]
}
